$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.766.14'
$ws.Range('E2').Value = '  +0.18%  '

$ws.Range('D3').Value = '2.620.41'
$ws.Range('E3').Value = '  -0.80%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.84'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.51%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.80'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.19%  '

$ws.Range('E7').Value = '  +0.02%  '

$ws.Range('E8').Value = '  -0.67%  '

$ws.Range('D9').Value = '2.619.85'
$ws.Range('E9').Value = '  -0.81%  '

$ws.Range('E10').Value = '  +6.18%  '

$ws.Range('E11').Value = '  -0.74%  '

$ws.Range('E12').Value = '  -0.54%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.347'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.42%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.54'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.84%  '

$ws.Range('E15').Value = '  +1.82%  '

$ws.Range('D16').Value = '3.097.26'
$ws.Range('E16').Value = '  -0.64%  '

$ws.Range('D17').Value = '67.662.76'
$ws.Range('E17').Value = '  +0.09%  '

$ws.Range('D18').Value = '2.618.93'
$ws.Range('E18').Value = '  -0.75%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '372.19'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.62%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.21'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.05%  '

$ws.Range('E21').Value = '  -2.68%  '

$ws.Range('E22').Value = '  -1.90%  '

$ws.Range('E23').Value = '  -3.36%  '

$ws.Range('E24').Value = '  -4.83%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '72.58'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +9.70%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.07%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.86'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.90%  '

$ws.Range('B28').Value = 'Bittensor'
$ws.Range('C28').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '595.40'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.03%  '

$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').Value = '2.757.57'
$ws.Range('E29').Value = '  -0.05%  '

$ws.Range('E30').Value = '  -1.53%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.30%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.81'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.71%  '

$ws.Range('E33').Value = '  -3.39%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.85'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.91%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.05%  '

$ws.Range('E36').Value = '  -3.58%  '

$ws.Range('E37').Value = '  -1.38%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '158.36'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.18%  '

$ws.Range('E39').Value = '  -2.05%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.89'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +2.75%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.368'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.27%  '

$ws.Range('E42').Value = '  -1.36%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.69'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +2.30%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '17.11'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +4.44%  '

$ws.Range('E45').Value = '  +0.06%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.42'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.82%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '156.23'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.19%  '

$ws.Range('D48').Value = '0.0₆0296'
$ws.Range('E48').Value = '  +2.08%  '

$ws.Range('E49').Value = '  -1.58%  '

$ws.Range('E50').Value = '  -2.85%  '

$ws.Range('E51').Value = '  -1.43%  '
